$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.098.53'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").Value = '1.652.37'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.37'
$ws.Range("E5").Value = '  +0.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5253'
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2596'
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06340'
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.36'
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.495'
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.663.97'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5479'
$ws.Range("E14").Value = '  +0.88%  '
$ws.Range("D15").Value = '0.0₅8236'
$ws.Range("E15").Value = '  +1.85%  '
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").Value = '26.097.94'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.572'
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '190.96'
$ws.Range("E20").Value = '  -0.72%  '
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.025'
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '142.11'
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1230'
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.227'
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.06'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.429'
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05849'
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.528'
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.256'
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.582'
$ws.Range("E33").Value = '  -0.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9472'
$ws.Range("E34").Value = '  -1.76%  '
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.779'
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5710'
$ws.Range("E38").Value = '  +1.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.762'
$ws.Range("E39").Value = '  -3.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8446'
$ws.Range("E40").Value = '  -1.61%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.39'
$ws.Range("E42").Value = '  +3.16%  '
$ws.Range("D43").Value = '1.024.66'
$ws.Range("E43").Value = '  +1.26%  '
$ws.Range("D44").Value = '1.797.35'
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.01'
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9997'
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4313'
$ws.Range("E47").Value = '  +3.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05146'
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.466'
$ws.Range("E49").Value = '  +1.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.795'
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09646'
$ws.Range("E51").Value = '  -0.41%  '
